# Fruta / hortaliza, semanal
# The weekly refresh shifts each week's record down one row (oldest record
# falls off row 8, newest becomes row 2) while keeping all the fields that
# never change (market/product descriptors) untouched. Only the date
# (D), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) columns move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot of the "before" values for the columns that change, keyed by row.
$data = @{
    2 = @{ D = 44232; M = 60; N = 11000; O = 12000; P = 11583; S = 827 }
    3 = @{ D = 44181; M = 65; N = 9000;  O = 10000; P = 9462;  S = 676 }
    4 = @{ D = 44172; M = 90; N = 8500;  O = 9000;  P = 8806;  S = 629 }
    5 = @{ D = 44210; M = 70; N = 10000; O = 11000; P = 10357; S = 740 }
    6 = @{ D = 44229; M = 55; N = 11000; O = 12000; P = 11364; S = 812 }
    7 = @{ D = 44216; M = 55; N = 11000; O = 12000; P = 11545; S = 825 }
    8 = @{ D = 44253; M = 90; N = 12000; O = 13000; P = 12667; S = 905 }
}

# New row -> old row that its values come from.
$mapping = @{
    2 = 8
    3 = 7
    4 = 6
    5 = 4
    6 = 5
    7 = 2
    8 = 3
}

foreach ($newRow in 2..8) {
    $oldRow = $mapping[$newRow]
    $vals = $data[$oldRow]

    $ws.Cells.Item($newRow, 4).Value  = $vals.D   # Column D - Fecha
    $ws.Cells.Item($newRow, 13).Value = $vals.M   # Column M - Volumen
    $ws.Cells.Item($newRow, 14).Value = $vals.N   # Column N - Precio minimo
    $ws.Cells.Item($newRow, 15).Value = $vals.O   # Column O - Precio maximo
    $ws.Cells.Item($newRow, 16).Value = $vals.P   # Column P - Precio promedio ponderado
    $ws.Cells.Item($newRow, 19).Value = $vals.S   # Column S - Precio $/Kg
}
